# Updated cryptos list on Wed Nov 13 12:51:07 UTC 2024 with GitHub Actions
# Refreshes the Price (D) and Volume(1h) (E) columns on Sheet1 for rows 2-51.
# Numeric-looking price strings are forced back to text (NumberFormat "@" then
# Style "Normal") so they keep matching the sheet's original inline-string
# representation instead of Excel auto-coercing them to floating-point numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '87.989.62'
$ws.Range("E2").Value = '  +1.31%  '
$ws.Range("D3").Value = '3.176.20'
$ws.Range("E3").Value = '  -2.11%  '
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '207.97'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.49%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '610.58'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.16%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.388'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +3.33%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.675'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +3.23%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.999'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.05%  '
$ws.Range("D10").Value = '3.172.11'
$ws.Range("E10").Value = '  -2.24%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.537'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -6.11%  '
$ws.Range("E12").Value = '  -1.05%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000245'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -5.00%  '
$ws.Range("B14").Value = 'Toncoin'
$ws.Range("C14").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.30'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.49%  '
$ws.Range("B15").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C15").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D15").Value = '3.748.65'
$ws.Range("E15").Value = '  -2.83%  '
$ws.Range("D16").Value = '87.672.33'
$ws.Range("E16").Value = '  +0.95%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '32.44'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -4.11%  '
$ws.Range("D18").Value = '3.177.73'
$ws.Range("E18").Value = '  -2.64%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.19'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +8.27%  '
$ws.Range("E20").Value = '  -3.38%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '412.18'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.75%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '8.49'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -4.74%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.07'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -4.52%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.26'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +3.03%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '12.30'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.85%  '
$ws.Range("D26").Value = '3.337.58'
$ws.Range("E26").Value = '  -3.13%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0000135'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +5.83%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '73.29'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -3.47%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.999'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.05%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.163'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -5.20%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.997'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.34%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '549.12'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.96%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '8.25'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -6.43%  '
$ws.Range("E34").Value = '  -6.38%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.86'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.98%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.86'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -4.83%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.131'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -4.11%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '21.89'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.64%  '
$ws.Range("E39").Value = '  +0.55%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.999'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.15%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.08'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +6.62%  '
$ws.Range("E42").Value = '  +0.05%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.92'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.75%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.373'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -4.99%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '152.37'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.31%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '173.97'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.80%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '43.23'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.39%  '
$ws.Range("E48").Value = '  +5.25%  '
$ws.Range("E49").Value = '  -7.01%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '24.17'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.63%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '3.96'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -6.13%  '
